$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 9 new rows to be inserted right after the header row (row 1),
# pushing the existing data rows (old rows 2-21) down to rows 11-30.
$insertedData = @(
    @(-0.05018257871270176, 0.03831957608461375, -0.02618168391287326),
    @(-0.0003359749913214444, 0.002840522676706245, 0.01902845310978596),
    @(0.0197737082280218, 0.003915645778179187, 0.0108062067255378),
    @(0.006835582219064174, 0.01575421430170528, -0.01602910399436943),
    @(0.009694431573152461, 0.04360967107117174, -0.02225992940366261),
    @(-0.006389650218188861, 0.1472857224941257, -0.09550878420472185),
    @(-0.0494189966470003, 0.3487124174833299, -0.1236696735024451),
    @(-0.1319224560260773, 0.6216225624084474, 0.00195476904511456),
    @(-0.2252138006687165, 0.7861163711547852, 0.1145494534075265)
)

# Insert 9 new rows before row 2 (shifting existing rows down)
$insertRange = $ws.Range("A2:C10")
$insertRange.EntireRow.Insert()

# Fill the newly inserted rows with their values
for ($i = 0; $i -lt $insertedData.Count; $i++) {
    $rowNum = 2 + $i
    $vals = $insertedData[$i]
    $ws.Cells.Item($rowNum, 1).Value = $vals[0]
    $ws.Cells.Item($rowNum, 2).Value = $vals[1]
    $ws.Cells.Item($rowNum, 3).Value = $vals[2]
}

# Append a new row (row 31) at the end with new data
$lastRow = 31
$ws.Cells.Item($lastRow, 1).Value = 0.002003637989982902
$ws.Cells.Item($lastRow, 2).Value = -0.01257160693407047
$ws.Cells.Item($lastRow, 3).Value = -0.01087340146303169
